$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '67.587.99'
$cell.Style = 'Normal'
$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  -3.17%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.634.44'
$cell.Style = 'Normal'
$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  -3.22%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.27%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '589.81'
$cell.Style = 'Normal'
$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  -2.14%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '184.94'
$cell.Style = 'Normal'
$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  -1.00%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.615'
$cell.Style = 'Normal'
$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -3.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +0.10%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.679'
$cell.Style = 'Normal'
$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -6.92%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  -10.87%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '54.08'
$cell.Style = 'Normal'
$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  -6.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  -14.09%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  -8.32%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '4.202.11'
$cell.Style = 'Normal'
$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -3.66%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '3.628.66'
$cell.Style = 'Normal'
$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -3.82%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  -0.45%  '
$cell.Style = 'Normal'
$cell = $ws.Range('B17')
$cell.NumberFormat = '@'
$cell.Value = 'WrappedBTC'
$cell.Style = 'Normal'
$cell = $ws.Range('C17')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$cell.Style = 'Normal'
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '67.367.06'
$cell.Style = 'Normal'
$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  -3.32%  '
$cell.Style = 'Normal'
$cell = $ws.Range('B18')
$cell.NumberFormat = '@'
$cell.Value = 'Chainlink'
$cell.Style = 'Normal'
$cell = $ws.Range('C18')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell.Style = 'Normal'
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '18.49'
$cell.Style = 'Normal'
$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  -5.97%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  -5.20%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -5.66%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '399.40'
$cell.Style = 'Normal'
$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -4.01%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  -7.13%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '85.95'
$cell.Style = 'Normal'
$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  -4.57%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -6.68%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '12.43'
$cell.Style = 'Normal'
$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -4.56%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  -0.47%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  -7.58%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -8.64%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '9.07'
$cell.Style = 'Normal'
$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  -5.60%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '31.43'
$cell.Style = 'Normal'
$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -5.53%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -8.61%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '66.74'
$cell.Style = 'Normal'
$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  +2.13%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '11.95'
$cell.Style = 'Normal'
$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  -5.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '595.58'
$cell.Style = 'Normal'
$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  -3.10%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -5.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -6.20%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  -0.05%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  -0.25%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  -7.63%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -18.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  -4.12%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  -9.47%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  -7.75%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '2.41'
$cell.Style = 'Normal'
$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  -13.19%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '2.714.48'
$cell.Style = 'Normal'
$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  -3.02%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  -4.24%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '3.06'
$cell.Style = 'Normal'
$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -5.75%  '
$cell.Style = 'Normal'
$cell = $ws.Range('B48')
$cell.NumberFormat = '@'
$cell.Value = 'Monero'
$cell.Style = 'Normal'
$cell = $ws.Range('C48')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.Style = 'Normal'
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '138.17'
$cell.Style = 'Normal'
$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  -2.97%  '
$cell.Style = 'Normal'
$cell = $ws.Range('B49')
$cell.NumberFormat = '@'
$cell.Value = 'WEMIXToken'
$cell.Style = 'Normal'
$cell = $ws.Range('C49')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell.Style = 'Normal'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '2.54'
$cell.Style = 'Normal'
$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -7.69%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '8.28'
$cell.Style = 'Normal'
$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  -12.02%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -7.06%  '
$cell.Style = 'Normal'
